# Applies the "Updated cryptos list" price/volume refresh described in the commit.
# Cells in columns B-E are stored as text (inlineStr) in the workbook; numeric-looking
# price values are written with a leading single-quote so Excel keeps them as text
# instead of silently converting them to numbers (which would drop formatting such as
# trailing zeros, e.g. "0.08820" -> 0.0882).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.358.81"
$ws.Range("E2").Value = "  -2.45%  "

# Row 3
$ws.Range("D3").Value = "1.987.60"
$ws.Range("E3").Value = "  -6.19%  "

# Row 4
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'329.68"
$ws.Range("E5").Value = "  -4.90%  "

# Row 6
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.4941"
$ws.Range("E7").Value = "  -4.87%  "

# Row 8
$ws.Range("D8").Value = "'0.4197"
$ws.Range("E8").Value = "  -6.22%  "

# Row 9
$ws.Range("D9").Value = "'51.78"
$ws.Range("E9").Value = "  -4.51%  "

# Row 10
$ws.Range("D10").Value = "'0.08820"
$ws.Range("E10").Value = "  -5.77%  "

# Row 11
$ws.Range("E11").Value = "  -5.56%  "

# Row 12
$ws.Range("E12").Value = "  -8.63%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.994.48"
$ws.Range("E13").Value = "  -5.44%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'8.014"
$ws.Range("E14").Value = "  -7.87%  "

# Row 15
$ws.Range("D15").Value = "'6.484"
$ws.Range("E15").Value = "  -7.03%  "

# Row 16
$ws.Range("D16").Value = "'96.16"
$ws.Range("E16").Value = "  -6.29%  "

# Row 17
$ws.Range("E17").Value = "  +0.06%  "

# Row 18
$ws.Range("E18").Value = "  -5.59%  "

# Row 19
$ws.Range("D19").Value = "'0.06642"
$ws.Range("E19").Value = "  -0.84%  "

# Row 20
$ws.Range("D20").Value = "'19.67"
$ws.Range("E20").Value = "  -8.92%  "

# Row 21
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("D22").Value = "'5.945"
$ws.Range("E22").Value = "  -5.63%  "

# Row 23
$ws.Range("D23").Value = "29.394.10"
$ws.Range("E23").Value = "  -2.43%  "

# Row 24
$ws.Range("E24").Value = "  -7.22%  "

# Row 25
$ws.Range("D25").Value = "'2.289"
$ws.Range("E25").Value = "  -1.73%  "

# Row 26
$ws.Range("D26").Value = "'157.29"
$ws.Range("E26").Value = "  -3.23%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.49"
$ws.Range("E27").Value = "  -7.59%  "

# Row 28
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'6.517"
$ws.Range("E28").Value = "  -2.87%  "

# Row 29
$ws.Range("D29").Value = "'2.339"
$ws.Range("E29").Value = "  -7.96%  "

# Row 30
$ws.Range("D30").Value = "'127.43"

# Row 31
$ws.Range("D31").Value = "'1.050"
$ws.Range("E31").Value = "  -9.34%  "

# Row 32
$ws.Range("D32").Value = "'0.09915"
$ws.Range("E32").Value = "  -6.21%  "

# Row 33
$ws.Range("D33").Value = "'1.561"
$ws.Range("E33").Value = "  -12.58%  "

# Row 34
$ws.Range("D34").Value = "'5.832"
$ws.Range("E34").Value = "  -7.22%  "

# Row 35
$ws.Range("D35").Value = "'3.779"
$ws.Range("E35").Value = "  -4.84%  "

# Row 36
$ws.Range("D36").Value = "'9.587"
$ws.Range("E36").Value = "  -10.89%  "

# Row 37
$ws.Range("D37").Value = "'0.02444"
$ws.Range("E37").Value = "  -7.28%  "

# Row 38
$ws.Range("D38").Value = "'0.06331"
$ws.Range("E38").Value = "  -7.84%  "

# Row 39
$ws.Range("D39").Value = "'1.281"
$ws.Range("E39").Value = "  -3.65%  "

# Row 40
$ws.Range("D40").Value = "'11.73"
$ws.Range("E40").Value = "  -7.87%  "

# Row 41
$ws.Range("D41").Value = "'0.6477"
$ws.Range("E41").Value = "  -9.24%  "

# Row 42
$ws.Range("D42").Value = "'0.2062"
$ws.Range("E42").Value = "  -8.40%  "

# Row 43
$ws.Range("E43").Value = "  +0.21%  "

# Row 44
$ws.Range("D44").Value = "'0.6313"
$ws.Range("E44").Value = "  -8.37%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.44"
$ws.Range("E45").Value = "  -8.49%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.202"
$ws.Range("E46").Value = "  -8.02%  "

# Row 47
$ws.Range("D47").Value = "'1.269"
$ws.Range("E47").Value = "  +0.53%  "

# Row 48
$ws.Range("D48").Value = "'3.530"
$ws.Range("E48").Value = "  -2.78%  "

# Row 49
$ws.Range("D49").Value = "'0.00000000331"
$ws.Range("E49").Value = "  -4.65%  "

# Row 50
$ws.Range("D50").Value = "'0.06982"
$ws.Range("E50").Value = "  -2.78%  "

# Row 51
$ws.Range("D51").Value = "'1.142"
$ws.Range("E51").Value = "  -5.29%  "
